$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.721.21"
$ws.Cells.Item(2, 5).Value = "  +1.43%  "

$ws.Cells.Item(3, 4).Value = "1.807.83"

$ws.Cells.Item(4, 5).Value = "  +0.35%  "

$ws.Cells.Item(5, 4).Value = "'328.04"
$ws.Cells.Item(5, 5).Value = "  -3.02%  "

$ws.Cells.Item(6, 4).Value = "'0.9993"
$ws.Cells.Item(6, 5).Value = "  +0.37%  "

$ws.Cells.Item(7, 4).Value = "'0.4382"
$ws.Cells.Item(7, 5).Value = "  +0.08%  "

$ws.Cells.Item(8, 4).Value = "'0.3770"
$ws.Cells.Item(8, 5).Value = "  +6.74%  "

$ws.Cells.Item(9, 4).Value = "'44.99"
$ws.Cells.Item(9, 5).Value = "  -1.48%  "

$ws.Cells.Item(10, 4).Value = "'0.07688"
$ws.Cells.Item(10, 5).Value = "  +2.97%  "

$ws.Cells.Item(11, 4).Value = "'1.140"
$ws.Cells.Item(11, 5).Value = "  -1.47%  "

$ws.Cells.Item(12, 4).Value = "'22.68"
$ws.Cells.Item(12, 5).Value = "  -1.37%  "

$ws.Cells.Item(13, 4).Value = "'1.001"
$ws.Cells.Item(13, 5).Value = "  +0.32%  "

$ws.Cells.Item(14, 4).Value = "'6.269"
$ws.Cells.Item(14, 5).Value = "  -0.44%  "

$ws.Cells.Item(15, 4).Value = "'7.515"
$ws.Cells.Item(15, 5).Value = "  +2.95%  "

$ws.Cells.Item(16, 4).Value = "1.806.10"
$ws.Cells.Item(16, 5).Value = "  -0.48%  "

$ws.Cells.Item(17, 4).Value = "'0.00001093"
$ws.Cells.Item(17, 5).Value = "  +0.40%  "

$ws.Cells.Item(18, 4).Value = "'0.06719"
$ws.Cells.Item(18, 5).Value = "  +0.74%  "

$ws.Cells.Item(19, 4).Value = "'81.09"
$ws.Cells.Item(19, 5).Value = "  -1.25%  "

$ws.Cells.Item(20, 4).Value = "'0.9997"
$ws.Cells.Item(20, 5).Value = "  +0.35%  "

$ws.Cells.Item(21, 4).Value = "'17.64"
$ws.Cells.Item(21, 5).Value = "  +1.69%  "

$ws.Cells.Item(22, 4).Value = "'6.284"
$ws.Cells.Item(22, 5).Value = "  -2.73%  "

$ws.Cells.Item(23, 4).Value = "28.709.20"
$ws.Cells.Item(23, 5).Value = "  +1.37%  "

$ws.Cells.Item(24, 4).Value = "'11.77"
$ws.Cells.Item(24, 5).Value = "  -2.74%  "

$ws.Cells.Item(25, 4).Value = "'2.449"
$ws.Cells.Item(25, 5).Value = "  +2.54%  "

$ws.Cells.Item(26, 4).Value = "'20.57"
$ws.Cells.Item(26, 5).Value = "  -0.97%  "

$ws.Cells.Item(27, 4).Value = "'154.75"
$ws.Cells.Item(27, 5).Value = "  -0.28%  "

$ws.Cells.Item(28, 4).Value = "'2.362"
$ws.Cells.Item(28, 5).Value = "  -4.38%  "

$ws.Cells.Item(29, 4).Value = "2.016.16"
$ws.Cells.Item(29, 5).Value = "  -0.28%  "

$ws.Cells.Item(30, 4).Value = "'1.308"
$ws.Cells.Item(30, 5).Value = "  -0.23%  "

$ws.Cells.Item(31, 4).Value = "'131.01"
$ws.Cells.Item(31, 5).Value = "  -1.44%  "

$ws.Cells.Item(32, 4).Value = "'3.967"
$ws.Cells.Item(32, 5).Value = "  -2.34%  "

$ws.Cells.Item(33, 4).Value = "'5.829"
$ws.Cells.Item(33, 5).Value = "  -2.58%  "

$ws.Cells.Item(34, 4).Value = "'0.09202"
$ws.Cells.Item(34, 5).Value = "  -1.31%  "

$ws.Cells.Item(35, 4).Value = "'0.2226"
$ws.Cells.Item(35, 5).Value = "  +2.64%  "

$ws.Cells.Item(36, 4).Value = "'12.20"
$ws.Cells.Item(36, 5).Value = "  -1.29%  "

$ws.Cells.Item(37, 4).Value = "'0.06322"
$ws.Cells.Item(37, 5).Value = "  +0.68%  "

$ws.Cells.Item(38, 4).Value = "'5.207"
$ws.Cells.Item(38, 5).Value = "  -0.32%  "

$ws.Cells.Item(39, 4).Value = "'0.6610"
$ws.Cells.Item(39, 5).Value = "  -3.28%  "

$ws.Cells.Item(40, 4).Value = "'0.02314"
$ws.Cells.Item(40, 5).Value = "  -2.91%  "

$ws.Cells.Item(41, 4).Value = "'1.205"
$ws.Cells.Item(41, 5).Value = "  -1.42%  "

$ws.Cells.Item(42, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(42, 4).Value = "'1.439"
$ws.Cells.Item(42, 5).Value = "  -3.68%  "

$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 4).Value = "'8.066"
$ws.Cells.Item(43, 5).Value = "  -2.39%  "

$ws.Cells.Item(44, 4).Value = "'0.9988"
$ws.Cells.Item(44, 5).Value = "  +0.35%  "

$ws.Cells.Item(45, 4).Value = "'14.02"
$ws.Cells.Item(45, 5).Value = "  +0.14%  "

$ws.Cells.Item(46, 4).Value = "'0.6078"
$ws.Cells.Item(46, 5).Value = "  -1.78%  "

$ws.Cells.Item(47, 4).Value = "'3.795"
$ws.Cells.Item(47, 5).Value = "  -1.91%  "

$ws.Cells.Item(48, 4).Value = "'127.86"
$ws.Cells.Item(48, 5).Value = "  -1.65%  "

$ws.Cells.Item(49, 4).Value = "'2.027"
$ws.Cells.Item(49, 5).Value = "  -1.26%  "

$ws.Cells.Item(50, 4).Value = "'0.07075"
$ws.Cells.Item(50, 5).Value = "  -0.71%  "

$ws.Cells.Item(51, 4).Value = "'1.145"
$ws.Cells.Item(51, 5).Value = "  -2.69%  "

